# Added Datadriven test to Maddies
$wb = $excel.ActiveWorkbook

# Rename Sheet4 -> Maddieslogindata
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Name = "Maddieslogindata"

$loginData = $wb.Worksheets.Item("LoginData")

# --- literal values (mirrors the LoginData layout: header + 3 data rows) ---
$ws4.Range("A1").Value = "LoginTestDataStart"
$ws4.Range("B1").Value = "UserName"
$ws4.Range("C1").Value = "Password"
$ws4.Range("D1").Value = "runMode"

$ws4.Range("B2").Value = "thimmaraju.g@winwire.com"
$ws4.Range("C2").Value = 123456
$ws4.Range("D2").Value = "Y"

$ws4.Range("B3").Value = "winwiretestinghyd@gmail.com"
$ws4.Range("C3").Value = 456789
$ws4.Range("D3").Value = "Y"

$ws4.Range("B4").Value = "cjsrsss@gmail.com"
$ws4.Range("C4").Value = 123456
$ws4.Range("D4").Value = "N"

# Hyperlinks on the email cells (added before the format paste below so the
# final cell style keeps the bordered hyperlink look used elsewhere)
$ws4.Hyperlinks.Add($ws4.Range("B2"), "mailto:thimmaraju.g@winwire.com") | Out-Null
$ws4.Hyperlinks.Add($ws4.Range("B3"), "mailto:winwiretestinghyd@gmail.com") | Out-Null
$ws4.Hyperlinks.Add($ws4.Range("B4"), "mailto:cjsrsss@gmail.com") | Out-Null

# --- formats copied from the LoginData sheet (keeps style indices aligned) ---
$loginData.Range("A1:D1").Copy() | Out-Null
$ws4.Range("A1:D1").PasteSpecial(-4122) | Out-Null

$loginData.Range("A2:D2").Copy() | Out-Null
$ws4.Range("A2:D2").PasteSpecial(-4122) | Out-Null

$loginData.Range("A3:D3").Copy() | Out-Null
$ws4.Range("A3:D3").PasteSpecial(-4122) | Out-Null

$loginData.Range("B4:D4").Copy() | Out-Null
$ws4.Range("B4:D4").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# Column widths to fit the new content (best-fit-like)
$ws4.Columns.Item(1).ColumnWidth = 16.666666666666668
$ws4.Columns.Item(2).ColumnWidth = 25.541666666666668
$ws4.Columns.Item(3).ColumnWidth = 8.166666666666666
$ws4.Columns.Item(4).ColumnWidth = 8.041666666666666

# Make Maddieslogindata the active/selected sheet (also drives workbook activeTab)
$ws4.Select()
$ws4.Range("C7").Select()

$wb.Save()
